$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 652
$ws.Range("I28").Value = 698
$ws.Range("K28").Value = 698
$ws.Range("M28").Value = -213
$ws.Range("H32").Value = 2319.3635
$ws.Range("J32").Value = 2481.3
$ws.Range("L32").Value = 2481.3
$ws.Range("N32").Value = -3133.3
$ws.Range("H64").Value = 6104.143
$ws.Range("I64").Value = 6593.6665
$ws.Range("J64").Value = 5737
$ws.Range("K64").Value = 6593.6665
$ws.Range("L64").Value = 5737
$ws.Range("M64").Value = -6345.6665
$ws.Range("N64").Value = -6233
$ws.Range("H67").Value = 6104.143
$ws.Range("I67").Value = 6593.6665
$ws.Range("J67").Value = 5737
$ws.Range("K67").Value = 6593.6665
$ws.Range("L67").Value = 5737
$ws.Range("M67").Value = -5735.6665
$ws.Range("N67").Value = -7453
$ws.Range("H74").Value = 9459
$ws.Range("I74").Value = 9188.5
$ws.Range("K74").Value = 9188.5
$ws.Range("M74").Value = -8252.5
$ws.Range("H76").Value = 12464
$ws.Range("I76").Value = 11135.2
$ws.Range("J76").Value = 13792.8
$ws.Range("K76").Value = 11135.2
$ws.Range("L76").Value = 13792.8
$ws.Range("M76").Value = -10820.2
$ws.Range("N76").Value = -14422.8
$ws.Range("H77").Value = 9459
$ws.Range("I77").Value = 9188.5
$ws.Range("K77").Value = 45942.5
$ws.Range("M77").Value = -41262.5
$ws.Range("H79").Value = 12464
$ws.Range("I79").Value = 11135.2
$ws.Range("J79").Value = 13792.8
$ws.Range("K79").Value = 11135.2
$ws.Range("L79").Value = 13792.8
$ws.Range("M79").Value = -10043.2
$ws.Range("N79").Value = -15976.8
$ws.Range("H87").Value = 178530.19
$ws.Range("J87").Value = 193383.2
$ws.Range("L87").Value = 193383.2
$ws.Range("N87").Value = -195879.2
$ws.Range("H88").Value = 2652.6
$ws.Range("I88").Value = 1749.5
$ws.Range("J88").Value = 2878.375
$ws.Range("K88").Value = 1749.5
$ws.Range("L88").Value = 2878.375
$ws.Range("M88").Value = -1343.5
$ws.Range("N88").Value = -3690.375
$ws.Range("H90").Value = 178530.19
$ws.Range("J90").Value = 193383.2
$ws.Range("L90").Value = 580149.6000000001
$ws.Range("N90").Value = -592629.6000000001
$ws.Range("H91").Value = 2652.6
$ws.Range("I91").Value = 1749.5
$ws.Range("J91").Value = 2878.375
$ws.Range("K91").Value = 1749.5
$ws.Range("L91").Value = 2878.375
$ws.Range("M91").Value = -345.5
$ws.Range("N91").Value = -5686.375
$ws.Range("H106").Value = 4978.1113
$ws.Range("J106").Value = 3999.3333
$ws.Range("L106").Value = 3999.3333
$ws.Range("N106").Value = -5261.3333
$ws.Range("H112").Value = 10349.048
$ws.Range("I112").Value = 14480
$ws.Range("K112").Value = 43440
$ws.Range("M112").Value = -42332
$ws.Range("H137").Value = 10834.805
$ws.Range("I137").Value = 2342.0833
$ws.Range("K137").Value = 7026.249899999999
$ws.Range("M137").Value = -4476.249899999999

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 240
$ws.Range("I5").Value = 80
$ws.Range("K5").Value = 80
$ws.Range("M5").Value = 32
$ws.Range("H51").Value = 46315.5
$ws.Range("J51").Value = 46315.5
$ws.Range("L51").Value = 46315.5
$ws.Range("N51").Value = -47827.5
$ws.Range("H61").Value = 1850821.1
$ws.Range("J61").Value = 2465814.5
$ws.Range("L61").Value = 2465814.5
$ws.Range("N61").Value = -2466238.5
$ws.Range("H63").Value = 2525.3333
$ws.Range("J63").Value = 2000
$ws.Range("L63").Value = 2000
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 2525.3333
$ws.Range("J66").Value = 2000
$ws.Range("L66").Value = 10000
$ws.Range("N66").Value = -16864
$ws.Range("H74").Value = 12256.637
$ws.Range("I74").Value = 3638.625
$ws.Range("J74").Value = 35238
$ws.Range("K74").Value = 3638.625
$ws.Range("L74").Value = 35238
$ws.Range("M74").Value = -2764.625
$ws.Range("N74").Value = -36986
$ws.Range("H77").Value = 12256.637
$ws.Range("I77").Value = 3638.625
$ws.Range("J77").Value = 35238
$ws.Range("K77").Value = 18193.125
$ws.Range("L77").Value = 176190
$ws.Range("M77").Value = -13825.125
$ws.Range("N77").Value = -184926
$ws.Range("H110").Value = 5428.091
$ws.Range("I110").Value = 5428.091
$ws.Range("K110").Value = 5428.091
$ws.Range("M110").Value = -3383.091
$ws.Range("H122").Value = 1540429.1
$ws.Range("I122").Value = 1668722.1
$ws.Range("K122").Value = 5006166.300000001
$ws.Range("M122").Value = -5003716.300000001
$ws.Range("H132").Value = 2314377
$ws.Range("I132").Value = 3600.0833
$ws.Range("K132").Value = 10800.2499
$ws.Range("M132").Value = -8270.249899999999
$ws.Range("H136").Value = 1850821.1
$ws.Range("J136").Value = 2465814.5
$ws.Range("L136").Value = 7397443.5
$ws.Range("N136").Value = -7402543.5
$ws.Range("H139").Value = 59999
$ws.Range("J139").Value = 59999
$ws.Range("L139").Value = 59999
$ws.Range("N139").Value = -70279

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 240
$ws.Range("I4").Value = 80
$ws.Range("K4").Value = 80
$ws.Range("M4").Value = 35
$ws.Range("H22").Value = 12707
$ws.Range("I22").Value = 12707
$ws.Range("K22").Value = 12707
$ws.Range("M22").Value = -12534
$ws.Range("H86").Value = 1168.2858
$ws.Range("I86").Value = 1110.3334
$ws.Range("K86").Value = 1110.3334
$ws.Range("M86").Value = 12.66660000000002
$ws.Range("H89").Value = 1168.2858
$ws.Range("I89").Value = 1110.3334
$ws.Range("K89").Value = 5551.666999999999
$ws.Range("M89").Value = 64.33300000000054
$ws.Range("H105").Value = 1575.2
$ws.Range("I105").Value = 1580.7142
$ws.Range("K105").Value = 1580.7142
$ws.Range("M105").Value = 166.2858000000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2422.4
$ws.Range("J122").Value = 1900
$ws.Range("L122").Value = 5700
$ws.Range("N122").Value = -10600

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2796952.5
$ws.Range("I4").Value = 2916322.8
$ws.Range("J4").Value = 2001150
$ws.Range("K4").Value = 8748968.399999999
$ws.Range("L4").Value = 6003450
$ws.Range("M4").Value = -8748856.399999999
$ws.Range("N4").Value = -6003674
$ws.Range("H63").Value = 23700
$ws.Range("I63").Value = 19500
$ws.Range("K63").Value = 58500
$ws.Range("M63").Value = -57751
$ws.Range("H66").Value = 23700
$ws.Range("I66").Value = 19500
$ws.Range("K66").Value = 175500
$ws.Range("M66").Value = -171756
$ws.Range("H68").Value = 1349.125
$ws.Range("J68").Value = 1349.125
$ws.Range("L68").Value = 4047.375
$ws.Range("N68").Value = -5669.375
$ws.Range("H71").Value = 1349.125
$ws.Range("J71").Value = 1349.125
$ws.Range("L71").Value = 12142.125
$ws.Range("N71").Value = -20254.125
$ws.Range("H87").Value = 1500
$ws.Range("I87").Value = 1500
$ws.Range("K87").Value = 4500
$ws.Range("M87").Value = -3252
$ws.Range("H90").Value = 1500
$ws.Range("I90").Value = 1500
$ws.Range("K90").Value = 13500
$ws.Range("M90").Value = -7260

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 10893.111
$ws.Range("I97").Value = 1148.5714
$ws.Range("J97").Value = 44999
$ws.Range("K97").Value = 1148.5714
$ws.Range("L97").Value = 44999
$ws.Range("M97").Value = -652.5714
$ws.Range("N97").Value = -45991
$ws.Range("H113").Value = 2478.6
$ws.Range("I113").Value = 2264.3333
$ws.Range("K113").Value = 2264.3333
$ws.Range("M113").Value = -94.33329999999978
$ws.Range("H132").Value = 746009.8
$ws.Range("I132").Value = 5575.643
$ws.Range("K132").Value = 16726.929
$ws.Range("M132").Value = -14196.929

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1981.6666
$ws.Range("I82").Value = 519.2
$ws.Range("K82").Value = 519.2
$ws.Range("M82").Value = -158.2
$ws.Range("H85").Value = 1981.6666
$ws.Range("I85").Value = 519.2
$ws.Range("K85").Value = 519.2
$ws.Range("M85").Value = 728.8

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 85000
$ws.Range("I56").Value = 85000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 85000
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -84286
$ws.Range("H60").Value = 85000
$ws.Range("J60").Value = 85000
$ws.Range("L60").Value = 85000
$ws.Range("N60").Value = -86644
$ws.Range("H76").Value = 80000
$ws.Range("J76").Value = 80000
$ws.Range("L76").Value = 80000
$ws.Range("N76").Value = -80630
$ws.Range("H79").Value = 80000
$ws.Range("J79").Value = 80000
$ws.Range("L79").Value = 80000
$ws.Range("N79").Value = -82184
$ws.Range("H126").Value = 7301.467
$ws.Range("I126").Value = 5410.909
$ws.Range("K126").Value = 16232.727
$ws.Range("M126").Value = -13762.727
$ws.Range("H136").Value = 445847.88
$ws.Range("I136").Value = 4663.385
$ws.Range("K136").Value = 13990.155
$ws.Range("M136").Value = -11440.155

Write-Output "edits applied"